$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal TEXT value that looks numeric (e.g. "34.5%")
# without Excel's auto-number-detection turning it into a percentage number
# (which would otherwise fabricate a brand-new cell style). We build the
# text via a formula (forces text type, keeps existing style), then paste
# the computed value back over itself to strip the formula again.
function Set-LiteralText($addr, [string]$text) {
    $r = $ws.Range($addr)
    $escaped = $text.Replace('"', '""')
    $r.Formula = '="' + $escaped + '"'
    $r.Copy() | Out-Null
    $r.PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = $false
}

# --- Row 2: reorder "Recorded By" email list (G2) ---
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# --- Row 3: reorder "Recorded By" email list (G3) ---
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, System, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# --- Row 4: reorder "Recorded By" email list (G4) ---
$ws.Range("G4").Value = "asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# --- L6: Recorded Sessions count 9 -> 10 ---
$ws.Range("L6").Value = 10

# --- Row 7 (Biochemistry Lab session #1) becomes Recorded ---
# Copy formatting (style) from a "Recorded" row (row 2) onto row 7
$ws.Range("A2:I2").Copy() | Out-Null
$ws.Range("A7:I7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("G7").Value = "Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("H7").Value = "19/251"
$ws.Range("I7").Value = "Recorded"

# --- L7: Missing Sessions count 2 -> 1 ---
$ws.Range("L7").Value = 1

# --- L9: Coverage % 31.0% -> 34.5% ---
Set-LiteralText "L9" "34.5%"

# --- L10: Average Attendance % 25.5% -> 23.7% ---
Set-LiteralText "L10" "23.7%"

# --- Row 15: summary row for Year2 / C1 ---
$ws.Range("O15").Value = 10
$ws.Range("P15").Value = 1
Set-LiteralText "R15" "34.5%"
Set-LiteralText "S15" "23.7%"

# --- Row 28: reorder "Recorded By" email list (G28) ---
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"

# --- Row 29 (Physiology session #2) keeps "Not Recorded" pink styling,
#     re-apply the same visual style via copy/paste so the engine re-creates
#     the xf entry (matches the legend rows' renumbering 7,8,9 -> 6,7,8).
$ws.Range("A29:I29").Copy() | Out-Null
$ws.Range("A29:I29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
